$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '29.476.10'
Set-TextValue 'E2' '  -1.20%  '
Set-TextValue 'D3' '1.852.16'
Set-TextValue 'E3' '  -0.28%  '
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '242.39'
Set-TextValue 'E5' '  -0.76%  '
Set-TextValue 'D6' '0.6310'
Set-TextValue 'E6' '  -3.63%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'D8' '0.07573'
Set-TextValue 'E8' '  +0.42%  '
Set-TextValue 'E9' '  -0.17%  '
Set-TextValue 'E10' '  -0.51%  '
Set-TextValue 'D11' '0.07713'
Set-TextValue 'E11' '  +1.00%  '
Set-TextValue 'D12' '1.927.42'
Set-TextValue 'E12' '  +3.44%  '
Set-TextValue 'D13' '5.011'
Set-TextValue 'E13' '  -0.99%  '
Set-TextValue 'D14' '0.6895'
Set-TextValue 'E14' '  +0.17%  '
Set-TextValue 'D15' '83.50'
Set-TextValue 'E15' '  -0.30%  '
Set-TextValue 'D16' '0.000009868'
Set-TextValue 'E16' '  +2.04%  '
Set-TextValue 'D17' '2.186.82'
Set-TextValue 'E17' '  +3.03%  '
Set-TextValue 'D18' '6.182'
Set-TextValue 'E18' '  +0.83%  '
Set-TextValue 'D19' '29.598.30'
Set-TextValue 'E19' '  -0.92%  '
Set-TextValue 'D20' '233.71'
Set-TextValue 'E20' '  -1.38%  '
Set-TextValue 'E21' '  -0.90%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D22' '1.002'
Set-TextValue 'E22' '  +0.18%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D23' '7.688'
Set-TextValue 'E23' '  -1.27%  '
Set-TextValue 'D24' '1.001'
Set-TextValue 'E24' '  -0.09%  '
Set-TextValue 'D25' '155.37'
Set-TextValue 'E25' '  -1.93%  '
Set-TextValue 'D26' '0.1398'
Set-TextValue 'E26' '  -2.65%  '
Set-TextValue 'E27' '  -1.04%  '
Set-TextValue 'E28' '  -0.85%  '
Set-TextValue 'E29' '  -1.06%  '
Set-TextValue 'D30' '0.05789'
Set-TextValue 'E30' '  -3.99%  '
Set-TextValue 'D31' '1.255'
Set-TextValue 'E31' '  -2.33%  '
Set-TextValue 'D32' '4.131'
Set-TextValue 'E32' '  -0.72%  '
Set-TextValue 'D33' '4.022'
Set-TextValue 'E33' '  -1.59%  '
Set-TextValue 'D34' '1.894'
Set-TextValue 'E34' '  +1.20%  '
Set-TextValue 'E35' '  -0.63%  '
Set-TextValue 'D36' '0.7244'
Set-TextValue 'E36' '  -0.70%  '
Set-TextValue 'D37' '2.592'
Set-TextValue 'E37' '  -0.61%  '
Set-TextValue 'D38' '1.255.83'
Set-TextValue 'E38' '  +4.26%  '
Set-TextValue 'D39' '2.803'
Set-TextValue 'E39' '  -0.35%  '
Set-TextValue 'E40' '  +0.65%  '
Set-TextValue 'D41' '0.9072'
Set-TextValue 'E41' '  -0.64%  '
Set-TextValue 'E42' '  -2.46%  '
Set-TextValue 'D43' '2.084.52'
Set-TextValue 'E43' '  +2.61%  '
Set-TextValue 'E44' '  +0.02%  '
Set-TextValue 'D45' '67.94'
Set-TextValue 'E45' '  +1.09%  '
Set-TextValue 'D46' '101.69'
Set-TextValue 'E46' '  +0.18%  '
Set-TextValue 'E47' '  +1.18%  '
Set-TextValue 'D48' '7.373'
Set-TextValue 'E48' '  +0.45%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '9.147'
Set-TextValue 'E49' '  -0.36%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D50' '0.4033'
Set-TextValue 'E50' '  -0.78%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D51' '1.713'
Set-TextValue 'E51' '  +1.88%  '
